$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 562, shifting existing rows 562:655 down to 565:658
$ws.Rows("562:564").Insert()

# Common column values shared by the three new rows (Palta / Hass, Macroferia Regional de Talca)
$A = 5
$B = "Macroferia Regional de Talca"
$C = "Maule"
$E = 7
$F = "Fruta"
$G = 100106
$H = "Oleaginosos"
$I = 100106002
$J = "Palta"
$K = "Hass"

# Row-specific data: row, fecha(D), calidad(L), volumen(M), precio(N/O/P/S), unidad(Q), origen(R), kgUnidad(T)
$newRows = @(
    @{ Row = 562; D = 44511; L = "Primera"; M = 150; N = 2500; Q = "`$/kilo (en caja de 17 kilos)"; R = "Cabildo"; T = 1 },
    @{ Row = 563; D = 44511; L = "Segunda"; M = 100; N = 2200; Q = "`$/kilo (en caja de 17 kilos)"; R = "Cabildo"; T = 1 },
    @{ Row = 564; D = 44511; L = "Tercera"; M = 60;  N = 2000; Q = "`$/kilo (en caja de 17 kilos)"; R = "Cabildo"; T = 1 }
)

foreach ($rd in $newRows) {
    $r = $rd.Row
    $ws.Cells.Item($r, 1).Value2 = $A
    $ws.Cells.Item($r, 2).Value2 = $B
    $ws.Cells.Item($r, 3).Value2 = $C
    $ws.Cells.Item($r, 4).Value2 = $rd.D
    $ws.Cells.Item($r, 5).Value2 = $E
    $ws.Cells.Item($r, 6).Value2 = $F
    $ws.Cells.Item($r, 7).Value2 = $G
    $ws.Cells.Item($r, 8).Value2 = $H
    $ws.Cells.Item($r, 9).Value2 = $I
    $ws.Cells.Item($r, 10).Value2 = $J
    $ws.Cells.Item($r, 11).Value2 = $K
    $ws.Cells.Item($r, 12).Value2 = $rd.L
    $ws.Cells.Item($r, 13).Value2 = $rd.M
    $ws.Cells.Item($r, 14).Value2 = $rd.N
    $ws.Cells.Item($r, 15).Value2 = $rd.N
    $ws.Cells.Item($r, 16).Value2 = $rd.N
    $ws.Cells.Item($r, 17).Value2 = $rd.Q
    $ws.Cells.Item($r, 18).Value2 = $rd.R
    $ws.Cells.Item($r, 19).Value2 = $rd.N
    $ws.Cells.Item($r, 20).Value2 = $rd.T
}
